# Update the fill test demo header row from Chinese to English labels.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A1").Value = "Date"
$ws.Range("B1").Value = "Number"
$ws.Range("C1").Value = "String 1"
$ws.Range("D1").Value = "String 2"
$ws.Range("E1").Value = "Image"

# B1 previously had no explicit style; align it with the rest of the header
# row (same font as A1/C1/D1/E1) by copying formats from A1.
$ws.Range("A1").Copy()
$ws.Range("B1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fix the selection so only E1 is active (was E1:E2 with E2 active).
$ws.Range("E1").Select()
